$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 17 data)
$ws.Range("D2").Value = 44874
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 2500

# Row 3 (was row 19 data)
$ws.Range("D3").Value = 44166
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = "`$/caja 18 kilos"
$ws.Range("R3").Value = "La Ligua"
$ws.Range("S3").Value = 667
$ws.Range("T3").Value = 18

# Row 4 (was row 18 data)
$ws.Range("D4").Value = 44921
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "`$/bandeja 7 kilos"
$ws.Range("S4").Value = 2143
$ws.Range("T4").Value = 7

# Row 6 (was row 2 data)
$ws.Range("D6").Value = 44511
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 28000
$ws.Range("O6").Value = 28000
$ws.Range("P6").Value = 28000
$ws.Range("Q6").Value = "`$/bandeja 10 kilos"
$ws.Range("R6").Value = "Provincia de Los Andes"
$ws.Range("S6").Value = 2800
$ws.Range("T6").Value = 10

# Row 7 (was row 3 data)
$ws.Range("D7").Value = 44511
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 3200
$ws.Range("O7").Value = 3200
$ws.Range("P7").Value = 3200
$ws.Range("S7").Value = 320

# Row 8 (was row 20 data)
$ws.Range("D8").Value = 44901
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 25000
$ws.Range("S8").Value = 2500

# Row 9 (was row 11 data)
$ws.Range("D9").Value = 44879
$ws.Range("M9").Value = 25
$ws.Range("N9").Value = 30000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 30000
$ws.Range("S9").Value = 3000

# Row 10 (was row 13 data)
$ws.Range("D10").Value = 44483
$ws.Range("M10").Value = 35
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("Q10").Value = "`$/bandeja 5 kilos"
$ws.Range("S10").Value = 2000
$ws.Range("T10").Value = 5

# Row 11 (was row 15 data)
$ws.Range("D11").Value = 44859
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = "`$/bandeja 5 kilos"
$ws.Range("S11").Value = 4000
$ws.Range("T11").Value = 5

# Row 12 (was row 9 data)
$ws.Range("D12").Value = 44503
$ws.Range("M12").Value = 50
$ws.Range("R12").Value = "Provincia de Quillota"

# Row 13 (was row 8 data)
$ws.Range("D13").Value = 44519
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 28000
$ws.Range("O13").Value = 28000
$ws.Range("P13").Value = 28000
$ws.Range("Q13").Value = "`$/bandeja 10 kilos"
$ws.Range("S13").Value = 2800
$ws.Range("T13").Value = 10

# Row 14 (was row 12 data)
$ws.Range("D14").Value = 44515
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 28000
$ws.Range("O14").Value = 28000
$ws.Range("P14").Value = 28000
$ws.Range("Q14").Value = "`$/bandeja 10 kilos"
$ws.Range("R14").Value = "Provincia de Los Andes"
$ws.Range("S14").Value = 2800
$ws.Range("T14").Value = 10

# Row 15 (was row 16 data)
$ws.Range("D15").Value = 44868
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("S15").Value = 2800

# Row 16 (was row 7 data)
$ws.Range("D16").Value = 44889
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 30000
$ws.Range("O16").Value = 30000
$ws.Range("P16").Value = 30000
$ws.Range("Q16").Value = "`$/bandeja 10 kilos"
$ws.Range("S16").Value = 3000
$ws.Range("T16").Value = 10

# Row 17 (was row 14 data)
$ws.Range("D17").Value = 44488
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = "`$/bandeja 5 kilos"
$ws.Range("R17").Value = "La Ligua"
$ws.Range("S17").Value = 2400
$ws.Range("T17").Value = 5

# Row 18 (was row 4 data)
$ws.Range("D18").Value = 44902
$ws.Range("M18").Value = 90
$ws.Range("N18").Value = 25000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 25000
$ws.Range("Q18").Value = "`$/bandeja 10 kilos"
$ws.Range("S18").Value = 2500
$ws.Range("T18").Value = 10

# Row 19 (was row 10 data)
$ws.Range("D19").Value = 44496
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 55
$ws.Range("N19").Value = 28000
$ws.Range("O19").Value = 28000
$ws.Range("P19").Value = 28000
$ws.Range("Q19").Value = "`$/bandeja 10 kilos"
$ws.Range("R19").Value = "Provincia de Quillota"
$ws.Range("S19").Value = 2800
$ws.Range("T19").Value = 10

# Row 20 (was row 6 data)
$ws.Range("D20").Value = 44858
$ws.Range("M20").Value = 90
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("Q20").Value = "`$/bandeja 5 kilos"
$ws.Range("S20").Value = 4000
$ws.Range("T20").Value = 5
